$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (dates are Excel serial numbers for 2021-04-22 .. 2021-04-26)
$data = @(
    @(44308, 18, 125, 310.3431153483291),
    @(44309, 22, 122, 302.8948805799692),
    @(44310, 15, 112, 278.0674313521029),
    @(44311, 2,  96,  238.3435125875168),
    @(44312, 13, 96,  238.3435125875168)
)

$startRow = 234
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Write values first
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Copy the date-column formatting (style s="2") from the row above it
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row[0]
}

$excel.CutCopyMode = $false
